# Update "paises.xlsx" (countries & provincias Spain COVID data snapshot)
# - Inserts "Singapur" ahead of "Israel" in the ranking (rows 27-29 shift)
# - Inserts "Armenia" ahead of "Croacia" in the ranking (rows 67-71 shift)
# - Updates the "Montenegro" row (row 125) active/recovered counts
# - Refreshes the "last updated" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $vals[$i]
    }
}

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 09:22"

# --- Singapur now ranks ahead of Israel (rows 27-29) -------------------
# Col order: Pais, Casos totales, Nuevos casos, Casos activos, Recuperados,
#            Casos criticos, Muertes hoy, Muertes
Set-Row 27 @("Singapur", 16169, 528, 1188, 14967, 22, 0, 14)
Set-Row 28 @("Israel",   15870,  36, 8412, 7239, 117, 4, 219)
Set-Row 29 @("Pakistan", 15759, 234, 4052, 11361, 111, 3, 346)

# --- Armenia now ranks ahead of Croacia (rows 67-71) --------------------
Set-Row 67 @("Armenia",    2066, 134,  929, 1105, 10, 2, 32)
Set-Row 68 @("Croacia",    2062,   0, 1288,  707, 19, 0, 67)
Set-Row 69 @("Uzbekistan", 2017,  15, 1096,  912,  8, 0,  9)
Set-Row 70 @("Irak",       2003,   0, 1346,  565,  0, 0, 92)
Set-Row 71 @("Afganistan", 1939,   0,  252, 1627,  7, 0, 60)

# --- Montenegro (row 125): updated activos / recuperados counts --------
$ws.Cells.Item(125, 4).Value = 206
$ws.Cells.Item(125, 5).Value = 109
